# Re-run of the power-flow case with the slack/ext-grid voltage set to
# 1.02 p.u. (was 1.05 p.u.) -- updates res_bus/vm_pu.xlsx (Sheet1) rows 2-25
# (time steps 0-23), columns B..F and I..N, per the "case with 380 kV done"
# commit. Column G (=1, slack bus itself) and column H (empty) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.057391338460802
$ws.Range("D2").Value = 1.063181363228568
$ws.Range("E2").Value = 1.053512598422141
$ws.Range("F2").Value = 1.071075376557564
$ws.Range("I2").Value = 1.042061356881811
$ws.Range("J2").Value = 1.062388005046422
$ws.Range("K2").Value = 1.065900697402111
$ws.Range("L2").Value = 1.056258361411837
$ws.Range("M2").Value = 1.073773515731665
$ws.Range("N2").Value = 1.063896717871195
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059000970849544
$ws.Range("D3").Value = 1.064682206650664
$ws.Range("E3").Value = 1.05490389588881
$ws.Range("F3").Value = 1.07274677526203
$ws.Range("I3").Value = 1.042437486610388
$ws.Range("J3").Value = 1.063646914115719
$ws.Range("K3").Value = 1.067214780274688
$ws.Range("L3").Value = 1.05746123813221
$ws.Range("M3").Value = 1.075259286954092
$ws.Range("N3").Value = 1.065157414735769
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060040198023495
$ws.Range("D4").Value = 1.065651419754685
$ws.Range("E4").Value = 1.055802333060478
$ws.Range("F4").Value = 1.073826563228899
$ws.Range("I4").Value = 1.042678012325772
$ws.Range("J4").Value = 1.064458835509638
$ws.Range("K4").Value = 1.068062641458161
$ws.Range("L4").Value = 1.058237230014548
$ws.Range("M4").Value = 1.076218487476548
$ws.Range("N4").Value = 1.065970489151197
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.060476546200029
$ws.Range("D5").Value = 1.066058423471067
$ws.Range("E5").Value = 1.056179607440711
$ws.Range("F5").Value = 1.074280105011825
$ws.Range("I5").Value = 1.042778448724185
$ws.Range("J5").Value = 1.06479953468371
$ws.Range("K5").Value = 1.068418507116492
$ws.Range("L5").Value = 1.058562902222986
$ws.Range("M5").Value = 1.076621219918624
$ws.Range("N5").Value = 1.066311672157175
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.06054977945117
$ws.Range("D6").Value = 1.066126734862315
$ws.Range("E6").Value = 1.056242928577723
$ws.Range("F6").Value = 1.074356233453975
$ws.Range("I6").Value = 1.04279527261156
$ws.Range("J6").Value = 1.064856702718864
$ws.Range("K6").Value = 1.068478225044539
$ws.Range("L6").Value = 1.058617551711726
$ws.Range("M6").Value = 1.076688810515823
$ws.Range("N6").Value = 1.066368921377496
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060046030647672
$ws.Range("D7").Value = 1.065656859927203
$ws.Range("E7").Value = 1.055807375891678
$ws.Range("F7").Value = 1.07383262503794
$ws.Range("I7").Value = 1.042679357032001
$ws.Range("J7").Value = 1.064463390422756
$ws.Range("K7").Value = 1.068067398800694
$ws.Range("L7").Value = 1.058241583831633
$ws.Range("M7").Value = 1.076223870819347
$ws.Range("N7").Value = 1.065975050532813
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.057935806580592
$ws.Range("D8").Value = 1.063688986126953
$ws.Range("E8").Value = 1.053983176923624
$ws.Range("F8").Value = 1.071640595736441
$ws.Range("I8").Value = 1.042189064705721
$ws.Range("J8").Value = 1.062814018551665
$ws.Range("K8").Value = 1.066345307896932
$ws.Range("L8").Value = 1.056665370463861
$ws.Range("M8").Value = 1.074276098852031
$ws.Range("N8").Value = 1.06432333636448
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.054199127487079
$ws.Range("D9").Value = 1.060206118832131
$ws.Range("E9").Value = 1.050754347831511
$ws.Range("F9").Value = 1.067764328128944
$ws.Range("I9").Value = 1.041303103694254
$ws.Range("J9").Value = 1.059886738520595
$ws.Range("K9").Value = 1.063291703701372
$ws.Range("L9").Value = 1.053869534480582
$ws.Range("M9").Value = 1.070826653164856
$ws.Range("N9").Value = 1.061391899260003
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.051695096387346
$ws.Range("D10").Value = 1.057873384561093
$ws.Range("E10").Value = 1.048591610589057
$ws.Range("F10").Value = 1.065170326149167
$ws.Range("I10").Value = 1.040697478459387
$ws.Range("J10").Value = 1.057920656372442
$ws.Range("K10").Value = 1.061242610292282
$ws.Range("L10").Value = 1.051992815934503
$ws.Range("M10").Value = 1.06851484080385
$ws.Range("N10").Value = 1.059423025049686
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050607608772679
$ws.Range("D11").Value = 1.056860585976367
$ws.Range("E11").Value = 1.047652586649429
$ws.Range("F11").Value = 1.064044617168927
$ws.Range("I11").Value = 1.040431638664477
$ws.Range("J11").Value = 1.057065749011765
$ws.Range("K11").Value = 1.060352040376714
$ws.Range("L11").Value = 1.051177022937083
$ws.Range("M11").Value = 1.067510776102008
$ws.Range("N11").Value = 1.058566903622498
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.050203169638549
$ws.Range("D12").Value = 1.05648396898901
$ws.Range("E12").Value = 1.047303398705877
$ws.Range("F12").Value = 1.06362609198489
$ws.Range("I12").Value = 1.040332349235973
$ws.Range("J12").Value = 1.056747650085097
$ws.Range("K12").Value = 1.060020736839528
$ws.Range("L12").Value = 1.050873516493597
$ws.Range("M12").Value = 1.067137354099525
$ws.Range("N12").Value = 1.058248352958867
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050289945946203
$ws.Range("D13").Value = 1.05656477373449
$ws.Range("E13").Value = 1.047378318640625
$ws.Range("F13").Value = 1.0637158847656
$ws.Range("I13").Value = 1.040353671865792
$ws.Range("J13").Value = 1.056815908361359
$ws.Range("K13").Value = 1.06009182565149
$ws.Range("L13").Value = 1.050938641740803
$ws.Range("M13").Value = 1.067217475739289
$ws.Range("N13").Value = 1.058316708169712
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.05057418793978
$ws.Range("D14").Value = 1.056829463302683
$ws.Range("E14").Value = 1.047623730735237
$ws.Range("F14").Value = 1.064010029710692
$ws.Range("I14").Value = 1.040423442503189
$ws.Range("J14").Value = 1.057039466099688
$ws.Range("K14").Value = 1.060324665111322
$ws.Range("L14").Value = 1.051151944941535
$ws.Range("M14").Value = 1.067479918531379
$ws.Range("N14").Value = 1.058540583385671
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050749252450841
$ws.Range("D15").Value = 1.056992491481459
$ws.Range("E15").Value = 1.047774884885394
$ws.Range("F15").Value = 1.064191210486487
$ws.Range("I15").Value = 1.040466358203412
$ws.Range("J15").Value = 1.057177134439999
$ws.Range("K15").Value = 1.060468057789749
$ws.Range("L15").Value = 1.05128330357942
$ws.Range("M15").Value = 1.067641555873751
$ws.Range("N15").Value = 1.058678447230815
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.051767200274957
$ws.Range("D16").Value = 1.05794054261098
$ws.Range("E16").Value = 1.048653876024288
$ws.Range("F16").Value = 1.065244982216542
$ws.Range("I16").Value = 1.040715045196583
$ws.Range("J16").Value = 1.057977317438832
$ws.Range("K16").Value = 1.06130164407827
$ws.Range("L16").Value = 1.05204689000101
$ws.Range("M16").Value = 1.068581412334663
$ws.Range("N16").Value = 1.05947976658129
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.052404859183113
$ws.Range("D17").Value = 1.058534496637523
$ws.Range("E17").Value = 1.049204555437132
$ws.Range("F17").Value = 1.065905310007958
$ws.Range("I17").Value = 1.040870073453024
$ws.Range("J17").Value = 1.058478285176326
$ws.Range("K17").Value = 1.06182364007138
$ws.Range("L17").Value = 1.052525014407823
$ws.Range("M17").Value = 1.069170138739254
$ws.Range("N17").Value = 1.059981445750433
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.052776484508767
$ws.Range("D18").Value = 1.058880679349615
$ws.Range("E18").Value = 1.049525512626469
$ws.Range("F18").Value = 1.066290228755655
$ws.Range("I18").Value = 1.040960151718563
$ws.Range("J18").Value = 1.058770146386245
$ws.Range("K18").Value = 1.062127793923362
$ws.Range("L18").Value = 1.05280359190348
$ws.Range("M18").Value = 1.069513240886062
$ws.Range("N18").Value = 1.060273721436747
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.052903146901447
$ws.Range("D19").Value = 1.05899867484699
$ws.Range("E19").Value = 1.049634909524082
$ws.Range("F19").Value = 1.066421435900349
$ws.Range("I19").Value = 1.040990807320916
$ws.Range("J19").Value = 1.058869605342623
$ws.Range("K19").Value = 1.062231448978373
$ws.Range("L19").Value = 1.052898528291482
$ws.Range("M19").Value = 1.069630180675642
$ws.Range("N19").Value = 1.06037332163625
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.052336476611014
$ws.Range("D20").Value = 1.058470798031908
$ws.Range("E20").Value = 1.049145498143024
$ws.Range("F20").Value = 1.065834487927743
$ws.Range("I20").Value = 1.040853476317421
$ws.Range("J20").Value = 1.058424571817428
$ws.Range("K20").Value = 1.061767667766707
$ws.Range("L20").Value = 1.052473747737096
$ws.Range("M20").Value = 1.069107004255947
$ws.Range("N20").Value = 1.059927656112404
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.050490499605989
$ws.Range("D21").Value = 1.056751530467313
$ws.Range("E21").Value = 1.047551473912764
$ws.Range("F21").Value = 1.063923422135286
$ws.Range("I21").Value = 1.040402911856574
$ws.Range("J21").Value = 1.056973649135857
$ws.Range("K21").Value = 1.060256113753447
$ws.Range("L21").Value = 1.05108914592051
$ws.Range("M21").Value = 1.067402648675065
$ws.Range("N21").Value = 1.058474672954202
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049326973763648
$ws.Range("D22").Value = 1.055668131646471
$ws.Range("E22").Value = 1.046546970691136
$ws.Range("F22").Value = 1.062719615189651
$ws.Range("I22").Value = 1.040116470746879
$ws.Range("J22").Value = 1.056058219171145
$ws.Range("K22").Value = 1.059302805606176
$ws.Range("L22").Value = 1.050215783211534
$ws.Range("M22").Value = 1.066328341001178
$ws.Range("N22").Value = 1.057557942973938
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049944058907176
$ws.Range("D23").Value = 1.056242695775238
$ws.Range("E23").Value = 1.04707969616969
$ws.Range("F23").Value = 1.063357993043789
$ws.Range("I23").Value = 1.040268618788717
$ws.Range("J23").Value = 1.056543810549967
$ws.Range("K23").Value = 1.059808453978576
$ws.Range("L23").Value = 1.05067903876889
$ws.Range("M23").Value = 1.066898112810744
$ws.Range("N23").Value = 1.058044223948216
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.052367376710764
$ws.Range("D24").Value = 1.058499581495292
$ws.Range("E24").Value = 1.049172184348631
$ws.Range("F24").Value = 1.065866490116716
$ws.Range("I24").Value = 1.040860976921132
$ws.Range("J24").Value = 1.058448843637703
$ws.Range("K24").Value = 1.061792960223244
$ws.Range("L24").Value = 1.05249691387943
$ws.Range("M24").Value = 1.069135532913692
$ws.Range("N24").Value = 1.059951962401448
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055167374222427
$ws.Range("D25").Value = 1.061108387619726
$ws.Range("E25").Value = 1.051590833358221
$ws.Range("F25").Value = 1.068768119048183
$ws.Range("I25").Value = 1.041534772835732
$ws.Range("J25").Value = 1.060646038336018
$ws.Range("K25").Value = 1.064083447794249
$ws.Range("L25").Value = 1.054594549682719
$ws.Range("M25").Value = 1.062152277368256
$ws.Range("N25").Value = 1.062152277368256
